$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.148787617683411
$ws.Range("B1").Value = 2.281025171279907
$ws.Range("C1").Value = 4.682409286499023
$ws.Range("D1").Value = 2.641624212265015
$ws.Range("E1").Value = 1.24669361114502
